# Apply the two changes captured in the commit:
#  1. The table on slide 5 gets a new table style (GUID) applied.
#  2. The presentation's theme color scheme is switched from the
#     "Integral" (Red Violet) palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 -------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{A0FF34A8-0245-4EFB-8A34-A52DC55732A4}")

# --- 2. Swap the theme colour scheme over to the Office palette -----------
# Order exposed by ThemeColorScheme.Item(n):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0x000000,
    0xFFFFFF,
    0x44546A,
    0xE7E6E6,
    0x5B9BD5,
    0xED7D31,
    0xA5A5A5,
    0xFFC000,
    0x4472C4,
    0x70AD47,
    0x0563C1,
    0x954F72
)

$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $hex = $officeColors[$i - 1]
    $r = ($hex -shr 16) -band 0xFF
    $g = ($hex -shr 8) -band 0xFF
    $b = $hex -band 0xFF
    # PowerPoint's RGB colour integers are stored little-endian (R + G*256 + B*65536)
    $themeColors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
